# metricas_recorrencia_mensal.xlsx
# "alterei modelagem de rfma recorrencia retencao e faturamentos da add para
#  banco padrao e carreguei dados novos para add e bibi"
#
# The sheet has columns:
#   A = yearmonth, B = retained_customers, C = prev_total_customers,
#   D = retention_rate (= B / C * 100, stored as a plain numeric value).
#
# Reloading the underlying data (new "banco padrao") shifted several of the
# monthly retained_customers / prev_total_customers counts by a handful of
# rows and recomputed the retention_rate accordingly. Apply the updated
# values for B, C and D on the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    [PSCustomObject]@{ Row = 2; B = $null; C = 51; D = 35.29411764705883 }
    [PSCustomObject]@{ Row = 4; B = $null; C = 44; D = 29.54545454545455 }
    [PSCustomObject]@{ Row = 6; B = $null; C = 94; D = 26.59574468085106 }
    [PSCustomObject]@{ Row = 8; B = $null; C = 93; D = 39.78494623655914 }
    [PSCustomObject]@{ Row = 9; B = 42; C = 87; D = 48.27586206896552 }
    [PSCustomObject]@{ Row = 10; B = 51; C = 90; D = 56.66666666666666 }
    [PSCustomObject]@{ Row = 11; B = $null; C = 120; D = 52.5 }
    [PSCustomObject]@{ Row = 12; B = 60; C = 98; D = 61.22448979591837 }
    [PSCustomObject]@{ Row = 13; B = 63; C = 121; D = 52.06611570247934 }
    [PSCustomObject]@{ Row = 14; B = 69; C = 127; D = 54.33070866141733 }
    [PSCustomObject]@{ Row = 15; B = 78; C = 138; D = 56.52173913043478 }
    [PSCustomObject]@{ Row = 16; B = 77; C = 126; D = 61.11111111111111 }
    [PSCustomObject]@{ Row = 17; B = 78; C = 116; D = 67.24137931034483 }
    [PSCustomObject]@{ Row = 18; B = 84; C = 142; D = 59.15492957746478 }
    [PSCustomObject]@{ Row = 19; B = 90; C = 142; D = 63.38028169014085 }
    [PSCustomObject]@{ Row = 20; B = 88; C = 144; D = 61.11111111111111 }
    [PSCustomObject]@{ Row = 21; B = 76; C = 138; D = 55.07246376811595 }
    [PSCustomObject]@{ Row = 22; B = 87; C = 111; D = 78.37837837837837 }
    [PSCustomObject]@{ Row = 24; B = 70; C = 103; D = 67.96116504854369 }
    [PSCustomObject]@{ Row = 25; B = $null; C = 107; D = 71.02803738317756 }
    [PSCustomObject]@{ Row = 26; B = $null; C = 125; D = 55.2 }
    [PSCustomObject]@{ Row = 27; B = 73; C = 97; D = 75.25773195876289 }
    [PSCustomObject]@{ Row = 28; B = $null; C = 122; D = 61.47540983606557 }
    [PSCustomObject]@{ Row = 31; B = $null; C = 120; D = 66.66666666666666 }
    [PSCustomObject]@{ Row = 33; B = 106; C = 140; D = 75.71428571428571 }
    [PSCustomObject]@{ Row = 34; B = $null; C = 141; D = 73.04964539007092 }
    [PSCustomObject]@{ Row = 35; B = $null; C = 139; D = 72.66187050359713 }
    [PSCustomObject]@{ Row = 36; B = $null; C = 142; D = 76.05633802816901 }
    [PSCustomObject]@{ Row = 37; B = 114; C = 159; D = 71.69811320754717 }
    [PSCustomObject]@{ Row = 38; B = 114; C = 154; D = 74.02597402597402 }
    [PSCustomObject]@{ Row = 39; B = $null; C = 162; D = 67.90123456790124 }
    [PSCustomObject]@{ Row = 41; B = 102; C = 141; D = 72.3404255319149 }
    [PSCustomObject]@{ Row = 42; B = 121; C = 177; D = 68.36158192090396 }
    [PSCustomObject]@{ Row = 43; B = 120; C = 159; D = 75.47169811320755 }
    [PSCustomObject]@{ Row = 44; B = 133; C = 185; D = 71.89189189189189 }
    [PSCustomObject]@{ Row = 45; B = 141; C = 194; D = 72.68041237113401 }
    [PSCustomObject]@{ Row = 46; B = 142; C = 207; D = 68.59903381642512 }
    [PSCustomObject]@{ Row = 47; B = 152; C = 214; D = 71.02803738317756 }
    [PSCustomObject]@{ Row = 48; B = 162; C = 223; D = 72.64573991031391 }
    [PSCustomObject]@{ Row = 49; B = 169; C = 224; D = 75.44642857142857 }
    [PSCustomObject]@{ Row = 50; B = 166; C = 223; D = 74.43946188340807 }
    [PSCustomObject]@{ Row = 51; B = 149; C = 236; D = 63.13559322033898 }
    [PSCustomObject]@{ Row = 52; B = 116; C = 194; D = 59.79381443298969 }
    [PSCustomObject]@{ Row = 53; B = 136; C = 186; D = 73.11827956989248 }
    [PSCustomObject]@{ Row = 54; B = 151; C = 206; D = 73.30097087378641 }
    [PSCustomObject]@{ Row = 55; B = 136; C = 223; D = 60.98654708520179 }
    [PSCustomObject]@{ Row = 56; B = 126; C = 184; D = 68.47826086956522 }
    [PSCustomObject]@{ Row = 57; B = 74; C = 208; D = 35.57692307692308 }
)

foreach ($u in $updates) {
    if ($null -ne $u.B) {
        $ws.Cells.Item($u.Row, 2).Value = $u.B
    }
    if ($null -ne $u.C) {
        $ws.Cells.Item($u.Row, 3).Value = $u.C
    }
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
}

Write-Host "Applied" $updates.Count "row updates to" $ws.Name
